$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added to the series. It belongs right
# after the existing row 15 (row 16), so insert a blank row there first,
# which pushes the former rows 16-116 down to 17-117.
$ws.Rows.Item(16).Insert()

# Fill in the new row with the new observation. All the "constant" columns
# (market/region/category/etc.) are identical to every other row in this
# sheet; only the date (D) and volume (J) are genuinely new values, the
# rest (K, L, M, P) happen to repeat the common values used throughout.
$ws.Cells.Item(16, 1).Value2 = 8
$ws.Cells.Item(16, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(16, 3).Value2 = "Coquimbo"
$ws.Cells.Item(16, 4).Value2 = 44473
$ws.Cells.Item(16, 5).Value2 = 4
$ws.Cells.Item(16, 6).Value2 = 100112037
$ws.Cells.Item(16, 7).Value2 = "Cebollín"
$ws.Cells.Item(16, 8).Value2 = "Sin especificar"
$ws.Cells.Item(16, 9).Value2 = "Primera"
$ws.Cells.Item(16, 10).Value2 = 2900
$ws.Cells.Item(16, 11).Value2 = 900
$ws.Cells.Item(16, 12).Value2 = 1000
$ws.Cells.Item(16, 13).Value2 = 950
$ws.Cells.Item(16, 14).Value2 = "`$/paquete 6 unidades"
$ws.Cells.Item(16, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(16, 16).Value2 = 158
$ws.Cells.Item(16, 17).Value2 = 6
$ws.Cells.Item(16, 18).Value2 = "Hortaliza"
